$d = $word.ActiveDocument

# --- Edit 1: drop the comma before "en een gedeeld ethos" in the opening paragraph ---
$old1 = @'
cultuur ondersteunen en in stand houden, en een gedeeld ethos
'@
$new1 = @'
cultuur ondersteunen en in stand houden en een gedeeld ethos
'@
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "edit1 found: $found1"

# --- Edit 2: rework the Berkshire/Schneider paragraph (added + reworded sentences) ---
$old2 = @'
Op dit moment zijn er honderden groepen die opkomen voor mensen die niet willen dat het het onderwijs religieus wordt, boeken worden verbonden en de toiletten gecontrolleerd worden en die niet willen dat het publieke onderwijs worden afgebroken. Zij zijn op zoek naar welke boodschappen het sterkste overkomen en op zoek naar die duidelijke visie op de toekomst waar elk onderwijs in elke gemeenschap dat onderwijs krijgt dat het nodig heeft.  Berkshire en Schneider helpen ons om de boodschap duidelijk te krijgen. Daar horen ouders, leerkrachten en studenten toe maar ook mensen uit de lokale gemeenschappen. Leerkrachten, vakbonden en politici alleen zijn niet voldoende. Die tegenkracht vraagt iets van ons allemaal. De kern van het vraagstuk is of we alleen maar interesse hebben in het onderwijs van onze eigen kinderen of in het onderwijs aan alle kinderen. Berskshire en  Is onderwijs individueel goed of publiek goed. Het is nodig om ons af te vragen waar de school voor is, welke idealen we nastreven, welke waarden, gewoonten en tradities. Democratie vraagt dat we over deze zaken blijven nadenken. Te lang hebben we gedacht dat het allemaal wel goed komt en ons te veel met bijzaken bezig gehouden. Deze tijd laat zien dat dat niet voldoende. Deze tijd vraagt dat we opkomen voor publeike onderwijs. Het publieke onderwijs is van de mensen, door de mensen, voor de mensen. Want, zoals iemand in dit boek schrijft: “
'@
$new2 = @'
Op dit moment zijn er honderden groepen die opkomen voor mensen die niet willen dat het het onderwijs religieus wordt, boeken verbonden en de toiletten gecontrolleerd worden en die niet willen dat het publieke onderwijs wordt afgebroken. Deze groepen zijn op zoek naar welke boodschappen het sterkste overkomen en op zoek naar die duidelijke visie op de toekomst waar elk onderwijs in elke gemeenschap dat onderwijs krijgt dat het nodig heeft. Daar hadden ze in de negentiende en zeker in de eerste helft van de twintigste eeuw uitgesproken ideeën. Daar besteden Berkshire en Schneider in deze gids en handleiding niet heel veel aandacht aan. Ze helpen ons met deze uitgave wel om de boodschap duidelijker te krijgen en daar zit de kracht in. Het is nodig om ons af te vragen waar de school voor is, welke idealen we nastreven, welke waarden, gewoonten en tradities. Democratie vraagt van ons dat we over deze zaken blijven nadenken. Niet alleen leerkrachten, vakbonden en politici moeten dat nu. Die tegenkracht die nu nodig is vraagt iets van ons allemaal omdat het de kern van onze democratie raakt. Het is nodig dat we niet alleen interesse hebben in het onderwijs van onze eigen kinderen maar in het onderwijs aan aan alle kinderen. Onderwijs is niet zozeer een individueel goed maar een publiek goed. Te lang is er gedacht dat het allemaal wel goed komt en is er teveel tijd besteed aan bijzaken. Deze tijd laat zien dat dat niet voldoende. Democratie vraagt dat we opkomen voor publeike onderwijs, deze tijd vraagt dat. Ook in onze eigen land, waar het allemaal niet zo hoogop wordt gespeeld maar elementen ervan wel. Er is veel te leren van dit Amerikaanse verhaal dat laat zien dat het publieke onderwijs van de mensen is, door de mensen, voor de mensen. Want, zoals iemand in dit boek schrijft: “
'@

$rng = $d.Content
$found2 = $rng.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "edit2 found: $found2"
if ($found2) {
    $rng.Text = $new2
}

Write-Output "done"
